# Auto-generated Excel COM-interop script
# Applies numeric "想去人数" (F column) updates and a row-order swap
# for two events (with updated title/location/time/links) across sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = 144

$ws.Range("F4").Value = 129

$ws.Range("F5").Value = 1916

$ws.Range("F7").Value = 3995

$ws.Range("F8").Value = 523

$ws.Range("F9").Value = 1038

$ws.Range("F10").Value = 1302

$ws.Range("F11").Value = 652

$ws.Range("F12").Value = 364

$ws.Range("F13").Value = 94

$ws.Range("F14").Value = 2166

$ws.Range("F15").Value = 389

$ws.Range("F16").Value = 647714

$ws.Range("F17").Value = 1597

$ws.Range("F18").Value = 471

$ws.Range("F19").Value = 1418

$ws.Range("F20").Value = 665

$ws.Range("F21").Value = 536

$ws.Range("F22").Value = 1242

$ws.Range("F23").Value = 2159

$ws.Range("F24").Value = 1104

$ws.Range("F25").Value = 2659

$ws.Range("F26").Value = 1526

$ws.Range("F27").Value = 752

$ws.Range("F28").Value = 1504

$ws.Range("F29").Value = 516

$ws.Range("F30").Value = 1070

$ws.Range("F31").Value = 1071

$ws.Range("F32").Value = 73

$ws.Range("F33").Value = 1993

$ws.Range("F34").Value = 1321

$ws.Range("F35").Value = 1196

$ws.Range("F36").Value = 1956

$ws.Range("F37").Value = 1126

$ws.Range("F41").Value = 2538

$ws.Range("F43").Value = 968

$ws.Range("F45").Value = 870

$ws.Range("F46").Value = 136

$ws.Range("F47").Value = 666

$ws.Range("F48").Value = 10

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("F2").Value = 63

# Row 9: event content updated
$ws.Range("C9").Value = "上海·【阿良良木健 领衔策划】小闹一场·唱聊会"
$ws.Range("D9").Value = "瑞虹路188号3楼 Modernsky Lab"
$ws.Range("E9").Value = "2024.07.12 22:00-07.12 23:59"
$ws.Range("F9").Value = 93
$ws.Range("G9").Value = 88
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=88313"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202406/XZ8xAGA41719566082242.jpeg"

# Row 10: event content updated
$ws.Range("C10").Value = "上海·夜鹿x夜游x真夜中   三夜0nly「夜³歌症候群」联合乐队现场"
$ws.Range("D10").Value = "虹许路731号4号楼 THE BOXX•城市乐园"
$ws.Range("E10").Value = "2024.07.12 19:30-07.12 22:00"
$ws.Range("F10").Value = 466
$ws.Range("G10").Value = 99
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85005"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202404/NPObaZdG1714384417870.png"

$ws.Range("F11").Value = 144324

$ws.Range("F12").Value = 144324

$ws.Range("F14").Value = 57

$ws.Range("F18").Value = 222

$ws.Range("F19").Value = 327

$ws.Range("F21").Value = 399

$ws.Range("F22").Value = 399

$ws.Range("F23").Value = 112

$ws.Range("F27").Value = 519

$ws.Range("F32").Value = 312

$ws.Range("F33").Value = 264

$ws.Range("F37").Value = 2

$ws.Range("F42").Value = 7

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("F4").Value = 3107

$ws.Range("F6").Value = 232

$ws.Range("F8").Value = 813

$ws.Range("F9").Value = 1137

$ws.Range("F10").Value = 629

$ws.Range("F11").Value = 1572

$ws.Range("F12").Value = 470

$ws.Range("F13").Value = 45

$ws.Range("F14").Value = 1798

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)

$ws.Range("F2").Value = 813

$ws.Range("F3").Value = 629

$ws.Range("F4").Value = 144

$ws.Range("F5").Value = 1572

$ws.Range("F6").Value = 470

$ws.Range("F7").Value = 129

$ws.Range("F8").Value = 1798

$ws.Range("F9").Value = 3995

$ws.Range("F11").Value = 523

$ws.Range("F12").Value = 1302

$ws.Range("F13").Value = 652

$ws.Range("F14").Value = 364

$ws.Range("F15").Value = 2166

$ws.Range("F17").Value = 389

$ws.Range("F18").Value = 647725

# Row 19: event content updated
$ws.Range("C19").Value = "上海·【阿良良木健 领衔策划】小闹一场·唱聊会"
$ws.Range("D19").Value = "瑞虹路188号3楼 Modernsky Lab"
$ws.Range("E19").Value = "2024.07.12 22:00-07.12 23:59"
$ws.Range("F19").Value = 93
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=88313"
$ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202406/XZ8xAGA41719566082242.jpeg"

# Row 20: event content updated
$ws.Range("C20").Value = "上海·夜鹿x夜游x真夜中   三夜0nly「夜³歌症候群」联合乐队现场"
$ws.Range("D20").Value = "虹许路731号4号楼 THE BOXX•城市乐园"
$ws.Range("E20").Value = "2024.07.12 19:30-07.12 22:00"
$ws.Range("F20").Value = 466
$ws.Range("G20").Value = 99
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=85005"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202404/NPObaZdG1714384417870.png"

$ws.Range("F21").Value = 1597

$ws.Range("F22").Value = 144324

$ws.Range("F23").Value = 471

$ws.Range("F24").Value = 1418

$ws.Range("F25").Value = 665

$ws.Range("F26").Value = 536

$ws.Range("F27").Value = 1242

$ws.Range("F28").Value = 2159

$ws.Range("F29").Value = 1104

$ws.Range("F30").Value = 2659

$ws.Range("F31").Value = 1526

$ws.Range("F33").Value = 1504

$ws.Range("F34").Value = 399

$ws.Range("F35").Value = 516

$ws.Range("F36").Value = 112

$ws.Range("F37").Value = 1071

$ws.Range("F38").Value = 1071

$ws.Range("F40").Value = 73

$ws.Range("F41").Value = 1993

$ws.Range("F42").Value = 1321

$ws.Range("F43").Value = 1196

$ws.Range("F44").Value = 1956

$ws.Range("F45").Value = 1126

$ws.Range("F46").Value = 312

$ws.Range("F47").Value = 312

$ws.Range("F48").Value = 2538

$ws.Range("F50").Value = 968

$ws.Range("F51").Value = 136

$ws.Range("F52").Value = 666

